# Add a new translation entry "SHOW_KEY_PRESS_TOTAL" to the KeyViewer sheet
# (adds the option to show the total number of key presses).
$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets("KeyViewer")

$ws.Range("A24").Value = "SHOW_KEY_PRESS_TOTAL"
$ws.Range("B24").Value = "Show key press total"
$ws.Range("C24").Value = "총 키를 누른 횟수 표시하기"
$ws.Range("D24").Value = "Mostrar número total de pulsaciones"
